$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header column (H1), matching the style of the other
# header cells in row 1 (e.g. G1 "sum") by copying its formatting over.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the corresponding data value for the new column (H2), a plain number.
$ws.Range("H2").Value = 0
